$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on D-column cells whose new value is a plain numeric-looking
# string (e.g. "92.60", "1.670") so Excel stores/display it verbatim instead of
# auto-converting to a Number and silently dropping significant trailing zeros.
$textCells = @(
    "D4", "D5", "D7", "D8", "D9", "D10", "D13", "D14", "D15", "D16", "D18", "D22", "D23", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values row by row, in sheet order.
$ws.Range("D2").Value = "26.923.38"
$ws.Range("E2").Value = "  +1.04%  "
$ws.Range("D3").Value = "1.819.19"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "309.73"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "0.4681"
$ws.Range("D8").Value = "0.3693"
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").Value = "0.07372"
$ws.Range("E9").Value = "  +1.82%  "
$ws.Range("D10").Value = "0.8715"
$ws.Range("E10").Value = "  +2.65%  "
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").Value = "1.799.68"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "5.359"
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("D14").Value = "92.60"
$ws.Range("E14").Value = "  +2.41%  "
$ws.Range("D15").Value = "0.07071"
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("D16").Value = "6.511"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "0.000008718"
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("E20").Value = "  +1.32%  "
$ws.Range("D21").Value = "26.961.02"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").Value = "5.342"
$ws.Range("E22").Value = "  +1.87%  "
$ws.Range("D23").Value = "10.57"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("D24").Value = "2.050.16"
$ws.Range("E24").Value = "  +1.85%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "151.36"
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("D27").Value = "2.192"
$ws.Range("E27").Value = "  +2.23%  "
$ws.Range("D28").Value = "18.40"
$ws.Range("E28").Value = "  +1.62%  "
$ws.Range("D29").Value = "5.325"
$ws.Range("E29").Value = "  +2.96%  "
$ws.Range("D30").Value = "115.83"
$ws.Range("E30").Value = "  +1.81%  "
$ws.Range("D31").Value = "0.08934"
$ws.Range("E31").Value = "  +1.25%  "
$ws.Range("D32").Value = "0.7699"
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("D33").Value = "1.165"
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("D34").Value = "4.501"
$ws.Range("E34").Value = "  +1.74%  "
$ws.Range("D35").Value = "2.903"
$ws.Range("E35").Value = "  +1.25%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").Value = "1.087"
$ws.Range("E37").Value = "  -2.09%  "
$ws.Range("E38").Value = "  +1.56%  "
$ws.Range("D39").Value = "0.05287"
$ws.Range("E39").Value = "  +1.80%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.955"
$ws.Range("E40").Value = "  +3.33%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "7.319"
$ws.Range("E41").Value = "  +2.94%  "
$ws.Range("D42").Value = "0.5347"
$ws.Range("E42").Value = "  +3.21%  "
$ws.Range("D43").Value = "2.366"
$ws.Range("E43").Value = "  +2.13%  "
$ws.Range("E44").Value = "  +1.68%  "
$ws.Range("D45").Value = "8.438"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").Value = "0.4959"
$ws.Range("E46").Value = "  +0.80%  "
$ws.Range("D47").Value = "10.45"
$ws.Range("E47").Value = "  +2.19%  "
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "1.001"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "1.670"
$ws.Range("E49").Value = "  +2.00%  "
$ws.Range("D50").Value = "103.38"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").Value = "0.06286"
$ws.Range("E51").Value = "  +0.21%  "
